$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Add the two new log rows (5 and 6) for the "Pattern Gen 1" task
# ---------------------------------------------------------------------

# Time values (From / To) for the new rows
$ws.Range("B5").Value = 0.34722222222222227
$ws.Range("C5").Value = 0.3611111111111111
$ws.Range("B6").Value = 0.375
$ws.Range("C6").Value = 0.3888888888888889
$ws.Range("B5:D6").NumberFormat = "h:mm"

# Duration formulas
$ws.Range("D5").Formula = "=C5-B5"
$ws.Range("D6").Formula = "=C6-B6"

# Date column - stored as text, matching the other rows ("2.3.2020" etc.)
$ws.Range("A5:A6").NumberFormat = "@"
$ws.Range("A5").Value = "3.3.2020"
$ws.Range("A6").Value = "3.3.2020"
$ws.Range("A5:A6").Style = "Normal"

# Task / Notes text
$ws.Range("E5").Value = "Pattern Gen 1"
$ws.Range("F5").Value = "Concept of module"

# Merge the Task/Notes cells across the two rows of this entry
$ws.Range("E5:E6").Merge() | Out-Null
$ws.Range("F5:F6").Merge() | Out-Null

# ---------------------------------------------------------------------
# Formatting
# ---------------------------------------------------------------------

# Every data/header cell is vertically centered
$ws.Range("A1:F6").VerticalAlignment = -4108

# Header row: bold
$ws.Range("A1:F1").Font.Bold = $true

# Task / Notes columns (E:F): left aligned
$ws.Range("E2:F6").HorizontalAlignment = -4131

# Date column (A): right aligned
$ws.Range("A2:A6").HorizontalAlignment = -4152

# Time columns (B:D): right aligned
$ws.Range("B2:D6").HorizontalAlignment = -4152

# Header row: centered horizontally
$ws.Range("A1:F1").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# View / selection
# ---------------------------------------------------------------------
$ws.Range("E11").Select() | Out-Null

# ---------------------------------------------------------------------
# Page setup
# ---------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
